$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.381.45"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.33%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.445.45"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.50%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.01%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'554.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -2.38%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'160.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -1.43%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.01%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.500"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -2.03%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.444.92"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -1.49%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -6.42%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -1.02%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.332"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -5.71%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  -1.98%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'2.895.01"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.50%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'68.232.52"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -1.36%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.0000166"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -4.13%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'23.20"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -3.38%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.450.38"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -2.02%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'10.66"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -4.05%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'339.42"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -1.72%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'6.95"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -5.20%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -2.61%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.07%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  -2.27%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  -4.42%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  -5.82%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'2.575.55"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -1.75%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'1.01"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +0.61%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'8.01"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -6.47%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.0₃0804"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -6.44%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'7.08"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -5.81%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'433.76"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -0.21%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'1.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.03%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -5.40%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -5.66%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'157.01"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +0.12%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -0.21%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +0.03%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  -2.56%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'17.76"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -1.57%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -3.47%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'4.38"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -3.46%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  -0.79%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'1.44"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -7.75%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'1.07"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +1.30%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  -5.65%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'131.01"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -4.63%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'3.31"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -2.51%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.0710"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -1.70%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.478"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -4.92%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -2.35%  "
$ws.Range("E51").Style = "Normal"
